$wb = $excel.ActiveWorkbook

# Sheet ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 746.75
$ws.Range("J33").Value = 1799.5
$ws.Range("L33").Value = 1799.5
$ws.Range("N33").Value = -2257.5

# Sheet ALC row 38
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1930.4
$ws.Range("J38").Value = 861.8
$ws.Range("L38").Value = 2585.4
$ws.Range("N38").Value = -3329.4

# Sheet ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2242.5
$ws.Range("I62").Value = 485
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 485
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = 139
$ws.Range("N62").Value = -5248

# Sheet ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 2242.5
$ws.Range("I65").Value = 485
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 2425
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = 695
$ws.Range("N65").Value = -26240

# Sheet ALC row 82
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 47000
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 47000
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 141000
$ws.Range("N82").Value = -141812
$ws.Range("M82").ClearContents()

# Sheet ALC row 85
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H85").Value = 47000
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 47000
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 141000
$ws.Range("N85").Value = -143808
$ws.Range("M85").ClearContents()

# Sheet ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 9999.25
$ws.Range("I116").Value = 9999.25
$ws.Range("K116").Value = 9999.25
$ws.Range("M116").Value = -6557.25

# Sheet ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1956.7142
$ws.Range("I132").Value = 2034.8334
$ws.Range("K132").Value = 6104.5002
$ws.Range("M132").Value = -3574.5002

# Sheet ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").ClearContents()

# Sheet ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3249.5
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# Sheet ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1512.9231
$ws.Range("I32").Value = 466.9
$ws.Range("K32").Value = 466.9
$ws.Range("M32").Value = -179.9

# Sheet ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3999
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 3999
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 3999
$ws.Range("N45").Value = -4753
$ws.Range("M45").ClearContents()

# Sheet ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1028
$ws.Range("I74").Value = 1028
$ws.Range("K74").Value = 1028
$ws.Range("M74").Value = -154

# Sheet ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1028
$ws.Range("I77").Value = 1028
$ws.Range("K77").Value = 5140
$ws.Range("M77").Value = -772

# Sheet ARM row 104
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H104").Value = 30332.666
$ws.Range("J104").Value = 30332.666
$ws.Range("L104").Value = 30332.666
$ws.Range("N104").Value = -37320.666

# Sheet ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 5844
$ws.Range("J122").Value = 4989.2
$ws.Range("L122").Value = 14967.6
$ws.Range("N122").Value = -19867.6

# Sheet BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3739.25
$ws.Range("I20").Value = 3564.8
$ws.Range("K20").Value = 3564.8
$ws.Range("M20").Value = -3317.8

# Sheet BSM row 88
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 26499.666
$ws.Range("J88").Value = 26499.666
$ws.Range("L88").Value = 26499.666
$ws.Range("N88").Value = -27311.666

# Sheet BSM row 91
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H91").Value = 26499.666
$ws.Range("J91").Value = 26499.666
$ws.Range("L91").Value = 26499.666
$ws.Range("N91").Value = -29307.666

# Sheet CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6298.8
$ws.Range("I58").Value = 1995
$ws.Range("J58").Value = 7374.75
$ws.Range("K58").Value = 1995
$ws.Range("L58").Value = 7374.75
$ws.Range("M58").Value = -1792
$ws.Range("N58").Value = -7780.75

# Sheet CRP row 59
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 140000
$ws.Range("J59").Value = 140000
$ws.Range("L59").Value = 140000
$ws.Range("N59").Value = -142290

# Sheet CRP row 88
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 23567.875
$ws.Range("J88").Value = 23567.875
$ws.Range("L88").Value = 23567.875
$ws.Range("N88").Value = -24379.875

# Sheet CRP row 91
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H91").Value = 23567.875
$ws.Range("J91").Value = 23567.875
$ws.Range("L91").Value = 23567.875
$ws.Range("N91").Value = -26375.875

# Sheet CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 6298.8
$ws.Range("I136").Value = 1995
$ws.Range("J136").Value = 7374.75
$ws.Range("K136").Value = 5985
$ws.Range("L136").Value = 22124.25
$ws.Range("M136").Value = -3435
$ws.Range("N136").Value = -27224.25

# Sheet CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 751.5
$ws.Range("J68").Value = 600
$ws.Range("L68").Value = 1800
$ws.Range("N68").Value = -3422

# Sheet CUL row 69
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 1658.8
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 1658.8
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 4976.4
$ws.Range("N69").Value = -6598.4
$ws.Range("M69").ClearContents()

# Sheet CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 751.5
$ws.Range("J71").Value = 600
$ws.Range("L71").Value = 5400
$ws.Range("N71").Value = -13512

# Sheet CUL row 72
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H72").Value = 1658.8
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 1658.8
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 14929.2
$ws.Range("N72").Value = -23041.2
$ws.Range("M72").ClearContents()

# Sheet GSM row 28
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H28").Value = 15000
$ws.Range("J28").Value = 15000
$ws.Range("L28").Value = 15000
$ws.Range("N28").Value = -15384

# Sheet GSM row 57
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 7800
$ws.Range("I57").Value = 7800
$ws.Range("K57").Value = 7800
$ws.Range("M57").Value = -6980

# Sheet GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

# Sheet GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2411.2
$ws.Range("I132").Value = 2761
$ws.Range("K132").Value = 8283
$ws.Range("M132").Value = -5753

# Sheet LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2687.5
$ws.Range("I16").Value = 3028.5715
$ws.Range("K16").Value = 3028.5715
$ws.Range("M16").Value = -2858.5715

# Sheet LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2000
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 2000
$ws.Range("N22").Value = -2590
$ws.Range("M22").ClearContents()

# Sheet LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2000
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 2000
$ws.Range("N27").Value = -2214
$ws.Range("M27").ClearContents()

# Sheet LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1996.5
$ws.Range("I46").Value = 1996.5
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 1996.5
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -1808.5
$ws.Range("N46").ClearContents()

# Sheet LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1084.25
$ws.Range("I68").Value = 993.5
$ws.Range("J68").Value = 1175
$ws.Range("K68").Value = 993.5
$ws.Range("L68").Value = 1175
$ws.Range("M68").Value = -244.5
$ws.Range("N68").Value = -2673

# Sheet LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1084.25
$ws.Range("I71").Value = 993.5
$ws.Range("J71").Value = 1175
$ws.Range("K71").Value = 4967.5
$ws.Range("L71").Value = 5875
$ws.Range("M71").Value = -1223.5
$ws.Range("N71").Value = -13363

# Sheet LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7900
$ws.Range("I122").Value = 7900
$ws.Range("K122").Value = 23700
$ws.Range("M122").Value = -21250

# Sheet WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1500
$ws.Range("I81").Value = 1500
$ws.Range("K81").Value = 3000
$ws.Range("M81").Value = -1939

# Sheet WVR row 82
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 47500.5
$ws.Range("J82").Value = 47500.5
$ws.Range("L82").Value = 47500.5
$ws.Range("N82").Value = -48266.5

# Sheet WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1500
$ws.Range("I84").Value = 1500
$ws.Range("K84").Value = 15000
$ws.Range("M84").Value = -9696

# Sheet WVR row 85
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H85").Value = 47500.5
$ws.Range("J85").Value = 47500.5
$ws.Range("L85").Value = 47500.5
$ws.Range("N85").Value = -50152.5

# Sheet WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3133.4
$ws.Range("I107").Value = 1585
$ws.Range("K107").Value = 4755
$ws.Range("M107").Value = -2835

# Sheet WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3095.2144
$ws.Range("I136").Value = 2564.1538
$ws.Range("J136").Value = 9999
$ws.Range("K136").Value = 7692.4614
$ws.Range("L136").Value = 29997
$ws.Range("M136").Value = -5142.4614
$ws.Range("N136").Value = -35097
